$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
# Row 41
$ws.Range("H41").Value = 5037.923
$ws.Range("I41").Value = 166.5
$ws.Range("J41").Value = 9213.429
$ws.Range("K41").Value = 166.5
$ws.Range("L41").Value = 9213.429
$ws.Range("M41").Value = 273.5
$ws.Range("N41").Value = -10093.429

# Row 76
$ws.Range("H76").Value = 74271.64
$ws.Range("I76").Value = 74271.64
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 74271.64
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = -73956.64
$ws.Range("N76").Value = $null

# Row 79
$ws.Range("H79").Value = 74271.64
$ws.Range("I79").Value = 74271.64
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 74271.64
$ws.Range("L79").Value = 0
$ws.Range("M79").Value = -73179.64
$ws.Range("N79").Value = $null

# Row 98
$ws.Range("H98").Value = 1339.4286
$ws.Range("I98").Value = 1332.16
$ws.Range("K98").Value = 1332.16
$ws.Range("M98").Value = 165.8399999999999

# Row 100
$ws.Range("H100").Value = 958.6667
$ws.Range("I100").Value = 918.1539
$ws.Range("K100").Value = 918.1539
$ws.Range("M100").Value = -377.1539

# Row 122
$ws.Range("H122").Value = 1339.4286
$ws.Range("I122").Value = 1332.16
$ws.Range("K122").Value = 3996.48
$ws.Range("M122").Value = -1546.48

# Row 125
$ws.Range("H125").Value = 2100
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 2100
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 18900
$ws.Range("M125").Value = $null
$ws.Range("N125").Value = -23820

# Row 132
$ws.Range("H132").Value = 5715623
$ws.Range("I132").Value = 6212373
$ws.Range("K132").Value = 18637119
$ws.Range("M132").Value = -18634589

# Row 137
$ws.Range("H137").Value = 1364.4445
$ws.Range("I137").Value = 1163.3334
$ws.Range("J137").Value = 1766.6666
$ws.Range("K137").Value = 3490.0002
$ws.Range("L137").Value = 5299.9998
$ws.Range("M137").Value = -940.0001999999999
$ws.Range("N137").Value = -10399.9998

$ws = $wb.Worksheets.Item(2)
# Row 32
$ws.Range("H32").Value = 20937.037
$ws.Range("I32").Value = 23364.146
$ws.Range("J32").Value = 4294
$ws.Range("K32").Value = 23364.146
$ws.Range("L32").Value = 4294
$ws.Range("M32").Value = -23077.146
$ws.Range("N32").Value = -4868

# Row 61
$ws.Range("H61").Value = 2435.4119
$ws.Range("I61").Value = 1954
$ws.Range("J61").Value = 4000
$ws.Range("K61").Value = 1954
$ws.Range("L61").Value = 4000
$ws.Range("M61").Value = -1742
$ws.Range("N61").Value = -4424

# Row 63
$ws.Range("H63").Value = 1430857.1
$ws.Range("I63").Value = 2002000
$ws.Range("J63").Value = 3000
$ws.Range("K63").Value = 2002000
$ws.Range("L63").Value = 3000
$ws.Range("M63").Value = -2001314
$ws.Range("N63").Value = -4372

# Row 66
$ws.Range("H66").Value = 1430857.1
$ws.Range("I66").Value = 2002000
$ws.Range("J66").Value = 3000
$ws.Range("K66").Value = 10010000
$ws.Range("L66").Value = 15000
$ws.Range("M66").Value = -10006568
$ws.Range("N66").Value = -21864

# Row 136
$ws.Range("H136").Value = 2435.4119
$ws.Range("I136").Value = 1954
$ws.Range("J136").Value = 4000
$ws.Range("K136").Value = 5862
$ws.Range("L136").Value = 12000
$ws.Range("M136").Value = -3312
$ws.Range("N136").Value = -17100

$ws = $wb.Worksheets.Item(3)
# Row 80
$ws.Range("H80").Value = 455.8125
$ws.Range("I80").Value = 722
$ws.Range("J80").Value = 334.81818
$ws.Range("K80").Value = 722
$ws.Range("L80").Value = 334.81818
$ws.Range("M80").Value = 276
$ws.Range("N80").Value = -2330.81818

# Row 83
$ws.Range("H83").Value = 455.8125
$ws.Range("I83").Value = 722
$ws.Range("J83").Value = 334.81818
$ws.Range("K83").Value = 3610
$ws.Range("L83").Value = 1674.0909
$ws.Range("M83").Value = 1382
$ws.Range("N83").Value = -11658.0909

# Row 107
$ws.Range("H107").Value = 868.2222
$ws.Range("I107").Value = 815
$ws.Range("K107").Value = 815
$ws.Range("M107").Value = 1105

$ws = $wb.Worksheets.Item(4)
# Row 31
$ws.Range("H31").Value = 3280578.8
$ws.Range("I31").Value = 1890.9302
$ws.Range("J31").Value = 11113000
$ws.Range("K31").Value = 1890.9302
$ws.Range("L31").Value = 11113000
$ws.Range("M31").Value = -1595.9302
$ws.Range("N31").Value = -11113590

# Row 34
$ws.Range("H34").Value = 3280578.8
$ws.Range("I34").Value = 1890.9302
$ws.Range("J34").Value = 11113000
$ws.Range("K34").Value = 1890.9302
$ws.Range("L34").Value = 11113000
$ws.Range("M34").Value = -1688.9302
$ws.Range("N34").Value = -11113404

# Row 58
$ws.Range("H58").Value = 943.13336
$ws.Range("I58").Value = 951.45
$ws.Range("K58").Value = 951.45
$ws.Range("M58").Value = -748.45

# Row 99
$ws.Range("H99").Value = 2055.111
$ws.Range("I99").Value = 1300
$ws.Range("J99").Value = 2999
$ws.Range("K99").Value = 1300
$ws.Range("L99").Value = 2999
$ws.Range("M99").Value = 198
$ws.Range("N99").Value = -5995

# Row 126
$ws.Range("H126").Value = 2055.111
$ws.Range("I126").Value = 1300
$ws.Range("J126").Value = 2999
$ws.Range("K126").Value = 3900
$ws.Range("L126").Value = 8997
$ws.Range("M126").Value = -1430
$ws.Range("N126").Value = -13937

# Row 136
$ws.Range("H136").Value = 943.13336
$ws.Range("I136").Value = 951.45
$ws.Range("K136").Value = 2854.35
$ws.Range("M136").Value = -304.3500000000004

$ws = $wb.Worksheets.Item(5)
# Row 82
$ws.Range("H82").Value = 3511.9092
$ws.Range("I82").Value = 1026.2
$ws.Range("J82").Value = 5583.3335
$ws.Range("K82").Value = 3078.6
$ws.Range("L82").Value = 16750.0005
$ws.Range("M82").Value = -2672.6
$ws.Range("N82").Value = -17562.0005

# Row 85
$ws.Range("H85").Value = 3511.9092
$ws.Range("I85").Value = 1026.2
$ws.Range("J85").Value = 5583.3335
$ws.Range("K85").Value = 3078.6
$ws.Range("L85").Value = 16750.0005
$ws.Range("M85").Value = -1674.6
$ws.Range("N85").Value = -19558.0005

# Row 107
$ws.Range("H107").Value = 334.07693
$ws.Range("J107").Value = 340.0909
$ws.Range("L107").Value = 1020.2727
$ws.Range("N107").Value = -4860.2727

# Row 122
$ws.Range("H122").Value = 927801.2
$ws.Range("I122").Value = 2315.2856
$ws.Range("J122").Value = 2223481.5
$ws.Range("K122").Value = 20837.5704
$ws.Range("L122").Value = 20011333.5
$ws.Range("M122").Value = -18387.5704
$ws.Range("N122").Value = -20016233.5

$ws = $wb.Worksheets.Item(6)
# Row 15
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").Value = $null

# Row 34
$ws.Range("H34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").Value = $null

# Row 74
$ws.Range("H74").Value = 59800
$ws.Range("J74").Value = 59800
$ws.Range("L74").Value = 59800
$ws.Range("N74").Value = -61672

# Row 76
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").Value = $null

# Row 77
$ws.Range("H77").Value = 59800
$ws.Range("J77").Value = 59800
$ws.Range("L77").Value = 179400
$ws.Range("N77").Value = -188760

# Row 79
$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").Value = $null

# Row 80
$ws.Range("H80").Value = 12500
$ws.Range("I80").Value = 3000
$ws.Range("J80").Value = 15666.667
$ws.Range("K80").Value = 3000
$ws.Range("L80").Value = 15666.667
$ws.Range("M80").Value = -2002
$ws.Range("N80").Value = -17662.667

# Row 81
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").Value = $null

# Row 83
$ws.Range("H83").Value = 12500
$ws.Range("I83").Value = 3000
$ws.Range("J83").Value = 15666.667
$ws.Range("K83").Value = 15000
$ws.Range("L83").Value = 78333.33499999999
$ws.Range("M83").Value = -10008
$ws.Range("N83").Value = -88317.33499999999

# Row 84
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").Value = $null

# Row 86
$ws.Range("H86").Value = 29743
$ws.Range("J86").Value = 29743
$ws.Range("L86").Value = 29743
$ws.Range("N86").Value = -32115

# Row 89
$ws.Range("H89").Value = 29743
$ws.Range("J89").Value = 29743
$ws.Range("L89").Value = 89229
$ws.Range("N89").Value = -101085

$ws = $wb.Worksheets.Item(7)
# Row 61
$ws.Range("H61").Value = 1267.6316
$ws.Range("J61").Value = 1981.25
$ws.Range("L61").Value = 1981.25
$ws.Range("N61").Value = -2385.25

# Row 93
$ws.Range("H93").Value = 1646.9333
$ws.Range("I93").Value = 2138.25
$ws.Range("J93").Value = 1085.4286
$ws.Range("K93").Value = 2138.25
$ws.Range("L93").Value = 1085.4286
$ws.Range("M93").Value = -890.25
$ws.Range("N93").Value = -3581.4286

# Row 113
$ws.Range("H113").Value = 1267.6316
$ws.Range("J113").Value = 1981.25
$ws.Range("L113").Value = 1981.25
$ws.Range("N113").Value = -6321.25

# Row 133
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").Value = $null

$ws = $wb.Worksheets.Item(8)
# Row 94
$ws.Range("H94").Value = 54749.5
$ws.Range("J94").Value = 54749.5
$ws.Range("L94").Value = 54749.5
$ws.Range("N94").Value = -56551.5

# Row 107
$ws.Range("H107").Value = 226
$ws.Range("I107").Value = 237.14285
$ws.Range("K107").Value = 711.4285500000001
$ws.Range("M107").Value = 1208.57145
